$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new column before column N (14th column)
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = 9.83

# Activate this worksheet (moves tabSelected here, and sets workbook activeTab)
$ws.Activate()

# Update the selection on this sheet to match the target (R7)
$ws.Range("R7").Select()
